# Updated remote access methods to use server_groups + units tests update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: GROUPS -> SERVER_GROUPS ---
$ws.Cells.Item(1, 8).Value = "SERVER_GROUPS"

# --- Remove the COMPLIANCE_GROUPS column (column I) entirely ---
$ws.Columns.Item(9).Delete()

# --- Row 2: HOST was blank, now the first server row ---
$ws.Cells.Item(2, 1).Value = "server01.example.com"

# --- Row 4: now a second SSH server entry (was a WinRM entry) ---
$ws.Cells.Item(4, 1).Value = "10.0.2.16"
$ws.Cells.Item(4, 2).Value = 22
$ws.Cells.Item(4, 3).Value = "CbwRam::RemoteAccess::Ssh::WithPassword"
$ws.Cells.Item(4, 4).Value = "master"
$ws.Cells.Item(4, 6).Value = "cyberwatch"
$ws.Cells.Item(4, 7).Value = "cyberwatch"
$ws.Cells.Item(4, 8).Value = "production, test"

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 14.9
$ws.Rows.Item(3).RowHeight = 1572.35

# --- Selection moves to row 2 ---
$ws.Rows.Item(2).Select()
